# "march 28 inputs but nine ball update failed"
# Appends the latest batch of 9-ball matchup rows (A:D = Player_1, Points_1,
# Player_2, Points_2) that hadn't made it into the sheet yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @(
    @(5,0,7,3),
    @(5,2,4,1),
    @(5,2,5,1),
    @(4,0,3,3),
    @(3,0,3,3),
    @(4,2,4,0),
    @(7,1,6,2),
    @(6,0,4,2),
    @(5,0,6,3),
    @(3,1,6,2),
    @(4,0,2,2),
    @(6,2,3,1),
    @(6,1,4,2),
    @(4,2,4,0),
    @(3,2,4,1),
    @(6,2,5,1),
    @(3,0,5,3),
    @(4,0,3,3),
    @(3,2,3,1),
    @(5,0,4,2),
    @(5,0,3,2),
    @(5,2,5,1)
)

# Last populated row before this edit.
$startRow = 1806

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + 1 + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

$lastRow = $startRow + $newRows.Count

# Mirror the author's cursor position after the paste: selection sitting one
# row below the newly-added data.
$target = $ws.Cells.Item($lastRow + 1, 1)
$target.Select()
